$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "toggle theme"
$ws.Range("C16").Value = "'26/08"
$ws.Range("E16").Value = "Lựu"
$ws.Range("F16").Value = "Đang làm"

$ws.Range("I18").Select() | Out-Null
